$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.921.74"
$ws.Range("E2").Value = "  -2.62%  "

$ws.Range("D3").Value = "3.468.28"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.38"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.20"
$ws.Range("E6").Value = "  -3.35%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -2.40%  "

$ws.Range("D9").Value = "3.464.38"
$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -5.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.85"
$ws.Range("E11").Value = "  -1.59%  "

$ws.Range("E12").Value = "  -4.08%  "

$ws.Range("D13").Value = "4.070.26"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.87"
$ws.Range("E15").Value = "  -6.68%  "

$ws.Range("D16").Value = "65.977.36"

$ws.Range("E17").Value = "  -2.72%  "

$ws.Range("D18").Value = "3.469.77"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.92"
$ws.Range("E19").Value = "  -3.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("E20").Value = "  -0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.75"
$ws.Range("E21").Value = "  -6.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.71"
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.23"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.533"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  +4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -6.86%  "

$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.77"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.98"
$ws.Range("E31").Value = "  -3.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  -5.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.11"
$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.28"

$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.27"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "29.16"
$ws.Range("E38").Value = "  +13.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.887"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("D40").Value = "2.807.54"
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("E41").Value = "  -5.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.55"
$ws.Range("E42").Value = "  -6.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  -3.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.43"
$ws.Range("E44").Value = "  -3.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0680"
$ws.Range("E45").Value = "  -4.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.97"
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.05"
$ws.Range("E47").Value = "  -7.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0288"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "308.94"
$ws.Range("E49").Value = "  -5.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.818"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.18"
$ws.Range("E51").Value = "  -2.02%  "
